# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.498.88"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "1.850.06"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.91"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.56"
$ws.Range("E8").Value = "  +9.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.308"
$ws.Range("E9").Value = "  +5.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0689"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +3.65%  "
$ws.Range("D12").Value = "2.117.56"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.62"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "1.841.75"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.73"
$ws.Range("E15").Value = "  +6.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.667"
$ws.Range("E16").Value = "  +5.17%  "
$ws.Range("D17").Value = "35.414.63"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.04"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.51"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "0.0₃0794"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  +7.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.83"
$ws.Range("E22").Value = "  +16.72%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.94"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.90"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").Value = "3.511.49"
$ws.Range("E29").Value = "  +44.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +7.15%  "
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.03"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0534"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.677"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "88.73"
$ws.Range("E37").Value = "  +9.40%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("E39").Value = "  +9.45%  "
$ws.Range("D40").Value = "1.339.17"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0195"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.42"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.28"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.92"
$ws.Range("E44").Value = "  +5.23%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").Value = "2.015.56"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.32"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("E51").Value = "  +0.15%  "
